# The workbook lists water-source statistics per Region/Province.
# A missing province - "NCR, CITY OF MANILA, FIRST DISTRICT (NOT A PROVINCE)" -
# needs to be added under the NCR region (rows were previously skipping it,
# which broke the downstream dashboard's graphs for that slice of data).
#
# It belongs right after "MIMAROPA REGION / ROMBLON" (row 17) and before the
# existing "NCR, FOURTH DISTRICT" row (old row 18), so we push rows 18..87
# down by one (columns B..I only - column A is a separate, untouched
# 0-based running index column) and populate the freed row 18 with the new
# province and zeroed counts (no data is available for it yet).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstShiftRow = 18
$lastDataRow = 87

# Walk bottom-up so we never overwrite a source row before it has been read.
for ($r = $lastDataRow; $r -ge $firstShiftRow; $r--) {
    for ($c = 2; $c -le 9; $c++) {
        $source = $ws.Cells.Item($r, $c)
        $target = $ws.Cells.Item($r + 1, $c)
        $target.Value = $source.Value2
    }
}

# Populate the newly vacated row with the missing province; region (column B)
# already carries over correctly from the shift above since row 18 itself
# was never overwritten (the loop only ever *writes* starting at row 19).
$ws.Cells.Item($firstShiftRow, 3).Value = "NCR, CITY OF MANILA, FIRST DISTRICT (NOT A PROVINCE)"
for ($c = 4; $c -le 9; $c++) {
    $ws.Cells.Item($firstShiftRow, $c).Value = 0
}

# Leave the selection on the cell that was actually edited.
$ws.Range("C18").Select() | Out-Null
